$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.033.82'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '2.350.03'
$ws.Range("E3").Value = '  -4.93%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '472.06'
$ws.Range("E5").Value = '  -3.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.39'
$ws.Range("E6").Value = '  -0.63%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("E8").Value = '  -2.16%  '
$ws.Range("D9").Value = '2.352.54'
$ws.Range("E9").Value = '  -5.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0963'
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.39'
$ws.Range("E11").Value = '  -6.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.319'
$ws.Range("E12").Value = '  -3.95%  '
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").Value = '2.757.60'
$ws.Range("E14").Value = '  -5.22%  '
$ws.Range("D15").Value = '54.932.83'
$ws.Range("E15").Value = '  -2.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.01'
$ws.Range("E16").Value = '  -5.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  -4.67%  '
$ws.Range("D18").Value = '2.350.52'
$ws.Range("E18").Value = '  -5.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.52'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '311.74'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.58'
$ws.Range("E21").Value = '  -5.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.59'
$ws.Range("E23").Value = '  -3.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '55.77'
$ws.Range("E24").Value = '  -4.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.391'
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.154'
$ws.Range("E27").Value = '  -6.44%  '
$ws.Range("D28").Value = '2.448.52'
$ws.Range("E28").Value = '  -5.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.16'
$ws.Range("E29").Value = '  -5.60%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '0.0₃0752'
$ws.Range("E31").Value = '  -4.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '148.33'
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.94'
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.02'
$ws.Range("E35").Value = '  -3.60%  '
$ws.Range("E36").Value = '  -5.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.54'
$ws.Range("E37").Value = '  -4.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.821'
$ws.Range("E38").Value = '  -4.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '33.51'
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.33'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.35'
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0945'
$ws.Range("E43").Value = '  +2.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0526'
$ws.Range("E44").Value = '  -5.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.575'
$ws.Range("E45").Value = '  -6.01%  '
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '254.94'
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("E48").Value = '  -3.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.43'
$ws.Range("E49").Value = '  -7.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.76'
$ws.Range("E50").Value = '  -4.87%  '
$ws.Range("D51").Value = '1.771.41'
$ws.Range("E51").Value = '  -6.51%  '
